$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.054.41'
$ws.Range("E2").Value = '  -0.56%  '

$ws.Range("D3").Value = '2.009.68'
$ws.Range("E3").Value = '  -1.63%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.24%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '225.77'
$ws.Range("E5").Value = '  -1.21%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.602'
$ws.Range("E6").Value = '  -1.56%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '55.20'
$ws.Range("E8").Value = '  -1.41%  '

$ws.Range("E9").Value = '  -2.70%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0776'
$ws.Range("E10").Value = '  -4.80%  '

$ws.Range("E11").Value = '  -4.43%  '

$ws.Range("D12").Value = '2.308.29'
$ws.Range("E12").Value = '  -1.48%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '13.98'
$ws.Range("E13").Value = '  -3.77%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.70'
$ws.Range("E14").Value = '  -3.91%  '

$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.734'
$ws.Range("E15").Value = '  -2.35%  '

$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.16'
$ws.Range("E16").Value = '  -1.92%  '

$ws.Range("D17").Value = '2.012.45'
$ws.Range("E17").Value = '  -1.50%  '

$ws.Range("D18").Value = '36.969.33'
$ws.Range("E18").Value = '  -0.58%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.15'
$ws.Range("E19").Value = '  +3.02%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '68.19'
$ws.Range("E20").Value = '  -2.06%  '

$ws.Range("D21").Value = '0.0₃0809'
$ws.Range("E21").Value = '  -4.14%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '223.72'
$ws.Range("E22").Value = '  -0.84%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  -0.03%  '

$ws.Range("E24").Value = '  +1.65%  '

$ws.Range("E25").Value = '  -4.81%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.06'
$ws.Range("E26").Value = '  -2.35%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.89'
$ws.Range("E27").Value = '  -6.35%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.56'
$ws.Range("E28").Value = '  -1.75%  '

$ws.Range("E29").Value = '  -3.57%  '

$ws.Range("E30").Value = '  -6.31%  '

$ws.Range("E31").Value = '  -1.57%  '

$ws.Range("E32").Value = '  -2.76%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0598'
$ws.Range("E33").Value = '  -1.88%  '

$ws.Range("E34").Value = '  -1.90%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.31'
$ws.Range("E35").Value = '  -2.85%  '

$ws.Range("E36").Value = '  +2.46%  '

$ws.Range("E37").Value = '  +0.18%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.11'
$ws.Range("E38").Value = '  -2.53%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.34'
$ws.Range("E39").Value = '  -0.89%  '

$ws.Range("D40").Value = '1.458.68'
$ws.Range("E40").Value = '  -2.55%  '

$ws.Range("E41").Value = '  -3.66%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '94.53'
$ws.Range("E42").Value = '  -0.94%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.27'
$ws.Range("E43").Value = '  +17.03%  '

$ws.Range("E44").Value = '  -3.02%  '

$ws.Range("E45").Value = '  -4.52%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '15.94'
$ws.Range("E46").Value = '  -4.63%  '

$ws.Range("E47").Value = '  -2.63%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.996'
$ws.Range("E48").Value = '  -1.45%  '

$ws.Range("E49").Value = '  -0.63%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.89'
$ws.Range("E50").Value = '  -0.43%  '

$ws.Range("D51").Value = '2.195.66'
$ws.Range("E51").Value = '  -1.43%  '

# Reset number format to default/general for the text-forced cells
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"